# INDIANA_2018.xlsx cleanup:
#  1. Rename header row columns to snake_case machine-readable names.
#  2. Title-case the Spanish linking particles ("de", "del", "la", "las",
#     "el", "los", "y") inside state/municipality names (e.g.
#     "Pabellón de Arteaga" -> "Pabellón De Arteaga").
#  3. Fix a floating point rounding blip in D258.
#  4. Drop the trailing footnote/source rows (1337:1341) that were appended
#     below the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 2. Title-case the linking particles across the whole used range -------
# Do this before touching the header row text so the header's own "de" particle
# (in e.g. "Estado de Origen") doesn't matter - it gets fully overwritten below.
$rng = $ws.UsedRange
$rng.Replace(" de ", " De ")   | Out-Null
$rng.Replace(" del ", " Del ") | Out-Null
$rng.Replace(" la ", " La ")   | Out-Null
$rng.Replace(" las ", " Las ") | Out-Null
$rng.Replace(" el ", " El ")   | Out-Null
$rng.Replace(" los ", " Los ") | Out-Null
$rng.Replace(" y ", " Y ")     | Out-Null

# --- 1. Rename header row to machine-readable snake_case names -------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 3. Correct the rounding artifact in D258 -------------------------------
$ws.Range("D258").Value = 0.09029168480626903

# --- 4. Remove the trailing footnote rows below the data table -------------
$ws.Rows("1337:1341").Delete() | Out-Null
